# Update column G ("K") values for rows 2-22 on Sheet1.
# These values are regenerated ("K instead of Strike#") per the commit message.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newK = @{
    2  = 7
    3  = 3
    4  = 3
    5  = 7
    6  = 2
    7  = 4
    8  = 5
    9  = 2
    10 = 5
    11 = 0
    12 = 3
    13 = 8
    14 = 5
    15 = 6
    16 = 5
    17 = 5
    18 = 2
    19 = 8
    20 = 3
    21 = 5
    22 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
